$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in assistant name cells
$ws.Range("C8").Value = "AndreiMarica(?)"
$ws.Range("D9").Value = "AndreiMarica"
$ws.Range("D10").Value = "AndreiMarica"
$ws.Range("D17").Value = "AndreiMarica"
$ws.Range("D18").Value = "AndreiMarica"
$ws.Range("D19").Value = "AndreiMarica"
$ws.Range("D20").Value = "AndreiMarica"

# Update view: scroll position, zoom, and selection
$win = $excel.ActiveWindow
$ws.Activate()
$win.ScrollRow = 5
$win.ScrollColumn = 1
$win.Zoom = 85
$ws.Range("C14").Select()
